$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- Step 1: Remove the "Meta description" paragraph (2nd paragraph). ---
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# --- Step 2: Insert a new bold title paragraph right before the final
#     "Prompt: ..." paragraph. We splice in two <w:p> fragments: the first
#     becomes the real new paragraph, the second (a lone empty run) is the
#     "merge forward" placeholder InsertXML needs so it doesn't spawn a
#     spurious blank paragraph of its own. ---
$count = $d.Paragraphs.Count
$promptPara = $d.Paragraphs($count)
$insertPoint = $d.Range($promptPara.Range.Start, $promptPara.Range.Start)
$insertPoint.InsertXML(
    '<w:p xmlns:w="' + $wNs + '"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Going Underground Slot for Free - Unique Features and Bonus Rounds</w:t></w:r></w:p>' +
    '<w:p xmlns:w="' + $wNs + '"><w:r/></w:p>')

# --- Step 3: Rebuild the (now shifted) prompt paragraph cleanly with the
#     replacement blurb text, preserving the italic run, and reusing the
#     same "merge forward" trick so no stray paragraph is left behind. ---
$count = $d.Paragraphs.Count
$promptPara = $d.Paragraphs($count)
$fullRange = $d.Range($promptPara.Range.Start, $promptPara.Range.End)
$fullRange.InsertXML(
    '<w:p xmlns:w="' + $wNs + '"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Going Underground slot with its unique modifiers and bonus rounds, and play for free on desktop or mobile.</w:t></w:r></w:p>' +
    '<w:p xmlns:w="' + $wNs + '"><w:r/></w:p>')

# That second InsertXML leaves one extra trailing empty paragraph (the
# placeholder had nothing left to merge into, since the prompt paragraph
# was the very last one in the body) - drop it.
$trailingCount = $d.Paragraphs.Count
$trailingPara = $d.Paragraphs($trailingCount)
$trailingPara.Range.Delete()
